$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (35 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1294.1111
$ws.Range("I9").Value = 341.5
$ws.Range("K9").Value = 341.5
$ws.Range("M9").Value = -172.5
$ws.Range("H17").Value = 41493.88
$ws.Range("J17").Value = 41493.88
$ws.Range("L17").Value = 124481.64
$ws.Range("N17").Value = -124817.64
$ws.Range("H62").Value = 3698.4
$ws.Range("I62").Value = 3698.4
$ws.Range("K62").Value = 3698.4
$ws.Range("M62").Value = -3074.4
$ws.Range("H65").Value = 3698.4
$ws.Range("I65").Value = 3698.4
$ws.Range("K65").Value = 18492
$ws.Range("M65").Value = -15372
$ws.Range("H92").Value = 3078.3044
$ws.Range("I92").Value = 1779
$ws.Range("K92").Value = 1779
$ws.Range("M92").Value = -531
$ws.Range("H137").Value = 3766
$ws.Range("I137").Value = 3741.4
$ws.Range("K137").Value = 11224.2
$ws.Range("M137").Value = -8674.200000000001
$ws.Range("H138").Value = 9240.209999999999
$ws.Range("I138").Value = 7368.4
$ws.Range("J138").Value = 9338.727000000001
$ws.Range("K138").Value = 22105.2
$ws.Range("L138").Value = 28016.181
$ws.Range("M138").Value = -16965.2
$ws.Range("N138").Value = -38296.181
$ws.Range("H141").Value = 3617.2856
$ws.Range("I141").Value = 3716.074
$ws.Range("K141").Value = 11148.222
$ws.Range("M141").Value = -5968.222

# --- Sheet: ARM (40 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15770.853
$ws.Range("I32").Value = 15770.853
$ws.Range("K32").Value = 15770.853
$ws.Range("M32").Value = -15483.853
$ws.Range("H45").Value = 2798.4138
$ws.Range("I45").Value = 2161.087
$ws.Range("J45").Value = 5241.5
$ws.Range("K45").Value = 2161.087
$ws.Range("L45").Value = 5241.5
$ws.Range("M45").Value = -1784.087
$ws.Range("N45").Value = -5995.5
$ws.Range("H61").Value = 6003557
$ws.Range("I61").Value = 8002951
$ws.Range("J61").Value = 1005071.4
$ws.Range("K61").Value = 8002951
$ws.Range("L61").Value = 1005071.4
$ws.Range("M61").Value = -8002739
$ws.Range("N61").Value = -1005495.4
$ws.Range("H110").Value = 8051.3076
$ws.Range("I110").Value = 8191.1113
$ws.Range("J110").Value = 7736.75
$ws.Range("K110").Value = 8191.1113
$ws.Range("L110").Value = 7736.75
$ws.Range("M110").Value = -6146.1113
$ws.Range("N110").Value = -11826.75
$ws.Range("H122").Value = 5443.467
$ws.Range("I122").Value = 5320.8276
$ws.Range("K122").Value = 15962.4828
$ws.Range("M122").Value = -13512.4828
$ws.Range("H132").Value = 8335724
$ws.Range("I132").Value = 2767.5
$ws.Range("K132").Value = 8302.5
$ws.Range("M132").Value = -5772.5
$ws.Range("H136").Value = 6003557
$ws.Range("I136").Value = 8002951
$ws.Range("J136").Value = 1005071.4
$ws.Range("K136").Value = 24008853
$ws.Range("L136").Value = 3015214.2
$ws.Range("M136").Value = -24006303
$ws.Range("N136").Value = -3020314.2

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3485.9524
$ws.Range("I94").Value = 3655.5264
$ws.Range("K94").Value = 3655.5264
$ws.Range("M94").Value = -3204.5264
$ws.Range("H134").Value = 4763964.5
$ws.Range("J134").Value = 20001504
$ws.Range("L134").Value = 60004512
$ws.Range("N134").Value = -60009582

# --- Sheet: CRP (18 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 20171.957
$ws.Range("I99").Value = 13581.333
$ws.Range("J99").Value = 27361.727
$ws.Range("K99").Value = 13581.333
$ws.Range("L99").Value = 27361.727
$ws.Range("M99").Value = -12083.333
$ws.Range("N99").Value = -30357.727
$ws.Range("H100").Value = 96031.164
$ws.Range("J100").Value = 96031.164
$ws.Range("L100").Value = 96031.164
$ws.Range("N100").Value = -98195.164
$ws.Range("H126").Value = 20171.957
$ws.Range("I126").Value = 13581.333
$ws.Range("J126").Value = 27361.727
$ws.Range("K126").Value = 40743.999
$ws.Range("L126").Value = 82085.181
$ws.Range("M126").Value = -38273.999
$ws.Range("N126").Value = -87025.181

# --- Sheet: CUL (4 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 4828
$ws.Range("I120").Value = 4828
$ws.Range("K120").Value = 14484
$ws.Range("M120").Value = -9646

# --- Sheet: GSM (32 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 222598.8
$ws.Range("I21").Value = 271498.5
$ws.Range("J21").Value = 27000
$ws.Range("K21").Value = 271498.5
$ws.Range("L21").Value = 27000
$ws.Range("M21").Value = -271325.5
$ws.Range("N21").Value = -27346
$ws.Range("H30").Value = 222598.8
$ws.Range("I30").Value = 271498.5
$ws.Range("J30").Value = 27000
$ws.Range("K30").Value = 271498.5
$ws.Range("L30").Value = 27000
$ws.Range("M30").Value = -271393.5
$ws.Range("N30").Value = -27210
$ws.Range("H102").Value = 1935.1666
$ws.Range("I102").Value = 1839.8
$ws.Range("J102").Value = 2412
$ws.Range("K102").Value = 1839.8
$ws.Range("L102").Value = 2412
$ws.Range("M102").Value = -217.8
$ws.Range("N102").Value = -5656
$ws.Range("H122").Value = 9222.579
$ws.Range("I122").Value = 6326.8125
$ws.Range("K122").Value = 18980.4375
$ws.Range("M122").Value = -16530.4375
$ws.Range("H126").Value = 46578930
$ws.Range("I126").Value = 81506130
$ws.Range("J126").Value = 9330.666999999999
$ws.Range("K126").Value = 244518390
$ws.Range("L126").Value = 27992.001
$ws.Range("M126").Value = -244515920
$ws.Range("N126").Value = -32932.001

# --- Sheet: LTW (59 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10234.467
$ws.Range("I7").Value = 9468
$ws.Range("J7").Value = 11110.429
$ws.Range("K7").Value = 9468
$ws.Range("L7").Value = 11110.429
$ws.Range("M7").Value = -9356
$ws.Range("N7").Value = -11334.429
$ws.Range("H22").Value = 2105.5
$ws.Range("I22").Value = 841.1429000000001
$ws.Range("K22").Value = 841.1429000000001
$ws.Range("M22").Value = -546.1429000000001
$ws.Range("H27").Value = 2105.5
$ws.Range("I27").Value = 841.1429000000001
$ws.Range("K27").Value = 841.1429000000001
$ws.Range("M27").Value = -734.1429000000001
$ws.Range("H46").Value = 1048.8334
$ws.Range("I46").Value = 1165
$ws.Range("J46").Value = 932.6667
$ws.Range("K46").Value = 1165
$ws.Range("L46").Value = 932.6667
$ws.Range("M46").Value = -977
$ws.Range("N46").Value = -1308.6667
$ws.Range("H55").Value = 1055.3062
$ws.Range("I55").Value = 800.06665
$ws.Range("J55").Value = 1458.3158
$ws.Range("K55").Value = 800.06665
$ws.Range("L55").Value = 1458.3158
$ws.Range("M55").Value = -627.06665
$ws.Range("N55").Value = -1804.3158
$ws.Range("H61").Value = 2455
$ws.Range("I61").Value = 2493.0303
$ws.Range("K61").Value = 2493.0303
$ws.Range("M61").Value = -2291.0303
$ws.Range("H93").Value = 18526552
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 18526552
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 18526552
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -18529048
$ws.Range("H113").Value = 2455
$ws.Range("I113").Value = 2493.0303
$ws.Range("K113").Value = 2493.0303
$ws.Range("M113").Value = -323.0302999999999
$ws.Range("H122").Value = 3267.9434
$ws.Range("I122").Value = 3288.255
$ws.Range("K122").Value = 9864.764999999999
$ws.Range("M122").Value = -7414.764999999999
$ws.Range("H126").Value = 10234.467
$ws.Range("I126").Value = 9468
$ws.Range("J126").Value = 11110.429
$ws.Range("K126").Value = 28404
$ws.Range("L126").Value = 33331.287
$ws.Range("M126").Value = -25934
$ws.Range("N126").Value = -38271.287
$ws.Range("H137").Value = 114617.7
$ws.Range("J137").Value = 117233.4
$ws.Range("L137").Value = 117233.4
$ws.Range("N137").Value = -127433.4

# --- Sheet: WVR (16 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2439.0667
$ws.Range("J100").Value = 3657.1667
$ws.Range("L100").Value = 7314.3334
$ws.Range("N100").Value = -8396.3334
$ws.Range("H103").Value = 60158
$ws.Range("J103").Value = 60158
$ws.Range("L103").Value = 60158
$ws.Range("N103").Value = -62502
$ws.Range("H122").Value = 2476.8096
$ws.Range("J122").Value = 3099.3333
$ws.Range("L122").Value = 9297.999899999999
$ws.Range("N122").Value = -14197.9999
$ws.Range("H136").Value = 291416.1
$ws.Range("I136").Value = 5705.125
$ws.Range("K136").Value = 17115.375
$ws.Range("M136").Value = -14565.375
